$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TL_ESS")
$ws.Range("E40").Value = 999
